# Regression logistic study workbook — "Correcting bugs for devtools check"
#  1. Refresh the "Saved on" timestamp in the intro sentence (A2).
#  2. Replace the tiny scientific-notation p-values (p=1e-28, p=3e-19, ...)
#     with the clearer "p<0.001" reporting convention, for every OR cell
#     whose p-value was below 0.001.
#  3. Right-align (instead of left-align) the numeric "OR (univariate/model1/
#     model2)" result cells in columns E:G so they line up with the rest of
#     the numeric table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Updated "Saved on" timestamp -------------------------------------
$ws.Range("A2").Value = "This is an interesting study. Saved on : 2022/09/28 17:56:24"

# --- 2. p=<tiny number> -> p<0.001 ----------------------------------------
$ws.Range("E5").Value  = "2.811 (2.372-3.42, p<0.001)"
$ws.Range("E6").Value  = "1.264 (1.203-1.333, p<0.001)"
$ws.Range("E9").Value  = "12.475 (6.162-28.139, p<0.001)"
$ws.Range("E10").Value = "65.932 (31.127-155.485, p<0.001)"

$ws.Range("F5").Value  = "4.043 (3.073-5.638, p<0.001)"
$ws.Range("F6").Value  = "1.331 (1.258-1.414, p<0.001)"
$ws.Range("F9").Value  = "12.686 (5.422-32.866, p<0.001)"
$ws.Range("F10").Value = "54.214 (20.658-157.555, p<0.001)"

$ws.Range("G5").Value  = "2.716 (2.204-3.46, p<0.001)"
$ws.Range("G6").Value  = "1.246 (1.162-1.341, p<0.001)"

# --- 3. Left -> right alignment for the OR value cells ---------------------
# xlRight = -4152 (standard Excel/VBA HorizontalAlignment enum value).
$xlRight = -4152
$orCells = @("E5", "F5", "G5", "E6", "F6", "G6", "E8", "F8", "E9", "F9", "E10", "F10", "G10")
foreach ($cellRef in $orCells) {
    $ws.Range($cellRef).HorizontalAlignment = $xlRight
}
